$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40: finish out the existing (previously mostly-empty) row with the new
# log entry - date / project / task / status.
$ws.Range("A40").Value = "4/8/2025(Remote)"
$ws.Range("B40").Value = "Car Tracking Project"
$ws.Range("C40").Value = "Add the code that shows the prompt before sending it to the model"
$ws.Range("F40").Value = "DONE"

# Row 41: brand-new log entry (date / project / task only).
$ws.Range("A41").Value = "4/8/2025(Remote)"
$ws.Range("B41").Value = "Car Tracking Project"
$ws.Range("C41").Value = "Check if there are any errors related to the json going in or out of the model"

# Update the sheet's visible selection to reflect where the author left off.
$ws.Range("D38").Select()
